$d = $word.ActiveDocument

# Locate the sentence that ends the "Submission Guidelines" bullet list and
# remove the trailing period, then append an ellipsis ("…") as a separate run,
# exactly as the source edit did.
$rng = $d.Content
$rng.Find.Execute("Provide a clear conclusion based on your analysis.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Text = "Provide a clear conclusion based on your analysis"
$rng.Collapse(0)

# Force a run boundary at the insertion point (without leaving any stray
# formatting behind) by bracketing it with a temporary bookmark.
$bm = $d.Bookmarks.Add("splitpoint", $rng)
$rng.InsertAfter([char]0x2026)
$d.Bookmarks("splitpoint").Delete()
